$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 100.7
$ws.Range("C7").Value = 101
$ws.Range("D7").Value = 99.59999999999999

$ws.Range("A8").Value = "2022年"
$ws.Range("B8").Value = 100.8
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Style = "Normal"

$ws.Range("A6").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
